# "Generate Report for Handoff" - a new handoff was generated for b.md in
# both the zh-cn and de-de locales. Update status + handoff file/datetime
# on the Overview, zh-cn and de-de sheets accordingly.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet: B3/C3 status for the b.md row ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# --- zh-cn sheet: row 3 (b.md) gets a fresh handoff file + datetime ---
$zhHandoffFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhHandoffDate = "2016-03-04 03:20:49"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $newStatus
$zhcn.Range("C3").Value = $zhHandoffFile
$zhcn.Range("D3").Value = $zhHandoffDate

foreach ($hl in $zhcn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$C$3') {
        $hl.TextToDisplay = $zhHandoffFile
    }
}

# --- de-de sheet: row 3 (b.md) gets a fresh handoff file + datetime ---
$deHandoffFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$deHandoffDate = "2016-03-04 03:21:03"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $newStatus
$dede.Range("C3").Value = $deHandoffFile
$dede.Range("D3").Value = $deHandoffDate

foreach ($hl in $dede.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$C$3') {
        $hl.TextToDisplay = $deHandoffFile
    }
}
